# Update NATMI TPM-derived ligand/receptor expression and specificity metrics
# (columns G-J, M-T for rows 2-9) to reflect results recomputed with the new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.27546433333333
$ws.Range("H2").Value = 48.826393
$ws.Range("I2").Value = 0.06628560529319844
$ws.Range("J2").Value = 0.06628560529319844
$ws.Range("M2").Value = 12.67919733333333
$ws.Range("N2").Value = 38.037592
$ws.Range("O2").Value = 0.9871416146107245
$ws.Range("P2").Value = 0.9871416146107247
$ws.Range("Q2").Value = 206.3598239739618
$ws.Range("R2").Value = 1857.238415765656
$ws.Range("S2").Value = 0.06543327943457709
$ws.Range("T2").Value = 0.06543327943457711
$ws.Range("G3").Value = 16.27546433333333
$ws.Range("H3").Value = 48.826393
$ws.Range("I3").Value = 0.06628560529319844
$ws.Range("J3").Value = 0.06628560529319844
$ws.Range("M3").Value = 0.1651576666666667
$ws.Range("N3").Value = 0.495473
$ws.Range("O3").Value = 0.01285838538927542
$ws.Range("P3").Value = 0.01285838538927542
$ws.Range("Q3").Value = 2.688017713209889
$ws.Range("R3").Value = 24.192159418889
$ws.Range("S3").Value = 0.0008523258586213401
$ws.Range("T3").Value = 0.0008523258586213401
$ws.Range("I4").Value = 0.3480686258826592
$ws.Range("J4").Value = 0.3480686258826592
$ws.Range("M4").Value = 12.67919733333333
$ws.Range("N4").Value = 38.037592
$ws.Range("O4").Value = 0.9871416146107245
$ws.Range("P4").Value = 0.9871416146107247
$ws.Range("Q4").Value = 1083.604502822191
$ws.Range("R4").Value = 9752.440525399721
$ws.Range("S4").Value = 0.3435930253491444
$ws.Range("T4").Value = 0.3435930253491445
$ws.Range("I5").Value = 0.3480686258826592
$ws.Range("J5").Value = 0.3480686258826592
$ws.Range("M5").Value = 0.1651576666666667
$ws.Range("N5").Value = 0.495473
$ws.Range("O5").Value = 0.01285838538927542
$ws.Range("P5").Value = 0.01285838538927542
$ws.Range("S5").Value = 0.004475600533514756
$ws.Range("T5").Value = 0.004475600533514756
$ws.Range("G6").Value = 42.61351133333333
$ws.Range("H6").Value = 127.840534
$ws.Range("I6").Value = 0.17355341356458
$ws.Range("J6").Value = 0.17355341356458
$ws.Range("M6").Value = 12.67919733333333
$ws.Range("N6").Value = 38.037592
$ws.Range("O6").Value = 0.9871416146107245
$ws.Range("P6").Value = 0.9871416146107247
$ws.Range("Q6").Value = 540.3051192615698
$ws.Range("R6").Value = 4862.746073354128
$ws.Range("S6").Value = 0.1713217968873423
$ws.Range("T6").Value = 0.1713217968873424
$ws.Range("G7").Value = 42.61351133333333
$ws.Range("H7").Value = 127.840534
$ws.Range("I7").Value = 0.17355341356458
$ws.Range("J7").Value = 0.17355341356458
$ws.Range("M7").Value = 0.1651576666666667
$ws.Range("N7").Value = 0.495473
$ws.Range("O7").Value = 0.01285838538927542
$ws.Range("P7").Value = 0.01285838538927542
$ws.Range("Q7").Value = 7.037948100286888
$ws.Range("R7").Value = 63.341532902582
$ws.Range("S7").Value = 0.002231616677237669
$ws.Range("T7").Value = 0.00223161667723767
$ws.Range("G8").Value = 101.183272
$ws.Range("H8").Value = 303.549816
$ws.Range("I8").Value = 0.4120923552595624
$ws.Range("J8").Value = 0.4120923552595624
$ws.Range("M8").Value = 12.67919733333333
$ws.Range("N8").Value = 38.037592
$ws.Range("O8").Value = 0.9871416146107245
$ws.Range("P8").Value = 0.9871416146107247
$ws.Range("Q8").Value = 1282.922672520341
$ws.Range("R8").Value = 11546.30405268307
$ws.Range("S8").Value = 0.4067935129396607
$ws.Range("T8").Value = 0.4067935129396608
$ws.Range("G9").Value = 101.183272
$ws.Range("H9").Value = 303.549816
$ws.Range("I9").Value = 0.4120923552595624
$ws.Range("J9").Value = 0.4120923552595624
$ws.Range("M9").Value = 0.1651576666666667
$ws.Range("N9").Value = 0.495473
$ws.Range("O9").Value = 0.01285838538927542
$ws.Range("P9").Value = 0.01285838538927542
$ws.Range("Q9").Value = 16.71119310921867
$ws.Range("R9").Value = 150.400737982968
$ws.Range("S9").Value = 0.005298842319901652
$ws.Range("T9").Value = 0.005298842319901652
